# Symbol-format refactor: append a new "SYMBOL_2017" data row (row 3) to
# Sheet1 of the VIN-upload refresh fixture, and move the active selection
# to D9 (matching the state the workbook was left in after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helpers --------------------------------------------------------------
# Columns B,C,D,E,F,H,I,J,O,R,S,X,Z,AA,AB,AD,AE,AF,AH,AI,AJ in row 2 all carry
# the existing "left aligned" cell style (style index 2 in styles.xml). Giving
# the brand-new row-3 cells that same HorizontalAlignment before writing their
# value makes the engine re-use that existing style entry (rather than minting
# a duplicate one), exactly like the rest of the data rows.
function Set-LeftCell($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.HorizontalAlignment = -4131   # xlLeft
    $cell.Value = $value
}

# Plain cells keep the worksheet's default (unstyled) formatting.
function Set-PlainCell($addr, $value) {
    $ws.Range($addr).Value = $value
}

# T3/V3 hold zero-padded numeric-looking codes ("0007"/"0002") that must be
# stored as literal text, not numbers. Prefixing with a leading apostrophe is
# exactly what Excel's UI does to force text entry, and it is what makes the
# engine mark the resulting style with quotePrefix="1" (new style indexes 3
# and 4) instead of silently coercing the value to a number.
function Set-LeftQuotedCell($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.HorizontalAlignment = -4131   # xlLeft
    $cell.Value = "'" + $value
}

function Set-QuotedCell($addr, $value) {
    $ws.Range($addr).Value = "'" + $value
}

# --- row 3 content (left to right, so new shared strings are interned in
# the same order they appear in the target file) --------------------------
Set-PlainCell     "A3"  "19XFB5F5&C"
Set-LeftCell       "B3"  "SYMBOL_2017"
Set-LeftCell       "C3"  2012
Set-LeftCell       "D3"  "HOND"
Set-LeftCell       "E3"  "HONDA"
Set-LeftCell       "F3"  "CIVIC"
Set-PlainCell     "G3"  "CIVIC NATURAL GAS"
Set-LeftCell       "H3"  26155
Set-LeftCell       "I3"  "4D SED"
Set-LeftCell       "J3"  "SEDAN 4 DOOR"
Set-PlainCell     "M3"  "SED"
Set-PlainCell     "N3"  "1.8L L4 COMPRESSED NATURAL GAS"
Set-LeftCell       "O3"  4
Set-LeftCell       "R3"  "2WD"
Set-LeftCell       "S3"  2
Set-LeftQuotedCell "T3"  "0007"
Set-PlainCell     "U3"  "FRNT/HEAD/SIDE/REAR SIDE AIRBAGS"
Set-QuotedCell     "V3"  "0002"
Set-PlainCell     "W3"  "4 WHEEL STANDARD"
Set-LeftCell       "X3"  "STD"
Set-PlainCell     "Y3"  "B-IMMOBILIZER/KEYLSS ENTRY/ALARM"
Set-LeftCell       "Z3"  51
Set-LeftCell       "AA3" 60
Set-LeftCell       "AB3" "Y"
Set-PlainCell     "AC3" "BI047"
Set-LeftCell       "AD3" "PD043"
Set-LeftCell       "AE3" "UM061"
Set-LeftCell       "AF3" "MP061"
Set-PlainCell     "AG3" 20180319
Set-LeftCell       "AH3" "Y"
Set-LeftCell       "AI3" "null"
Set-LeftCell       "AJ3" "null"

# --- selection -------------------------------------------------------------
$ws.Range("D9").Select()
